$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.845.42"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "2.276.77"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'249.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.08%  "
$ws.Range("D6").Value = "'0.642"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("D7").Value = "'79.16"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.70%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  +0.78%  "
$ws.Range("D10").Value = "'41.12"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.0973"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.56%  "
$ws.Range("D12").Value = "'7.35"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("E13").Value = "  +0.75%  "
$ws.Range("D14").Value = "2.620.05"
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("D15").Value = "'15.11"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.28%  "
$ws.Range("D16").Value = "'0.867"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.98%  "
$ws.Range("D17").Value = "2.286.31"
$ws.Range("E17").Value = "  -0.41%  "
$ws.Range("D18").Value = "42.796.77"
$ws.Range("D19").Value = "0.0₃0996"
$ws.Range("E19").Value = "  -2.18%  "
$ws.Range("D20").Value = "'6.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.12%  "
$ws.Range("D21").Value = "'72.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.94%  "
$ws.Range("D22").Value = "'233.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.44%  "
$ws.Range("D23").Value = "'2.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.03%  "
$ws.Range("D24").Value = "'3.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.68%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").Value = "'11.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.60%  "
$ws.Range("D27").Value = "'2.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.22%  "
$ws.Range("E28").Value = "  +2.21%  "
$ws.Range("D29").Value = "'168.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").Value = "'20.88"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.81%  "
$ws.Range("D31").Value = "'6.46"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.88%  "
$ws.Range("D32").Value = "'0.0856"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.69%  "
$ws.Range("E33").Value = "  -5.34%  "
$ws.Range("D34").Value = "'0.127"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.37%  "
$ws.Range("D35").Value = "'30.10"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.84%  "
$ws.Range("D36").Value = "'4.55"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.30%  "
$ws.Range("D37").Value = "'4.76"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.36%  "
$ws.Range("E38").Value = "  -2.38%  "
$ws.Range("D39").Value = "'13.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.90%  "
$ws.Range("D40").Value = "'2.26"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.81%  "
$ws.Range("D41").Value = "'5.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.73%  "
$ws.Range("D42").Value = "'0.209"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.82%  "
$ws.Range("D43").Value = "'111.81"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +16.39%  "
$ws.Range("D44").Value = "'61.28"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("D45").Value = "'8.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.05%  "
$ws.Range("E46").Value = "  -0.76%  "
$ws.Range("D47").Value = "'4.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.72%  "
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("E49").Value = "  -3.00%  "
$ws.Range("D50").Value = "'1.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.65%  "
$ws.Range("D51").Value = "'4.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.71%  "
